# Duplicate the last two data rows (row 9 "a1" and row 10 "b2") down into
# new rows 11 and 12, extending the used range from A1:FI10 to A1:FI12.
# Using Range.Copy(Destination) duplicates both the values and the
# formatting (e.g. the bold/centered/bordered style on column A) in one
# step, just like a manual copy/paste of the rows in Excel would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9:FI9").Copy($ws.Range("A11:FI11"))
$ws.Range("A10:FI10").Copy($ws.Range("A12:FI12"))
